# Generate Report for Handoff
# Updates the localization-status workbook:
#  - Overview / zh-cn / de-de sheets get two new handoff rows (3 and 4)
#    for a .md dependency file and a .png file, alongside the refreshed
#    values for the existing .png handoff row (row 2, which used to be
#    the 9d5b49ee...md handoff).

$wb = $excel.ActiveWorkbook

$M = [Type]::Missing

function Add-Hlink($ws, $cellRef, $url, $text) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, $M, $M, $text)
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 refresh: the old 9d5b49ee...md handoff is now the .png handoff
$ov.Range("D2").Value = "2016-43-14 08:43:03"
$ov.Range("A2").Value = "3de2e418-5a16-46cf-9df6-9bb7cac9905f.png"
Add-Hlink $ov "A2" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/3de2e418-5a16-46cf-9df6-9bb7cac9905f.png" "3de2e418-5a16-46cf-9df6-9bb7cac9905f.png"

# Row 3: new .md dependency row
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-43-14 08:43:03"
$ov.Range("A3").Value = "b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"
Add-Hlink $ov "A3" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/b1b1466b-8fdb-4993-b7a1-5da72f4af826.md" "b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"

# Row 4: new .png row
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-43-14 08:43:03"
$ov.Range("A4").Value = "e84bb446-0bb7-4b54-8a62-0c8b45e01771.png"
Add-Hlink $ov "A4" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/e84bb446-0bb7-4b54-8a62-0c8b45e01771.png" "e84bb446-0bb7-4b54-8a62-0c8b45e01771.png"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 refresh
$zh.Range("E2").Value = "2016-03-14 08:42:59"
$zh.Range("I2").Value = "IsDependency"
$zh.Range("J2").Value = "e2e\b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"

$zh.Range("A2").Value = "3de2e418-5a16-46cf-9df6-9bb7cac9905f.png"
Add-Hlink $zh "A2" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/3de2e418-5a16-46cf-9df6-9bb7cac9905f.png" "3de2e418-5a16-46cf-9df6-9bb7cac9905f.png"

$zh.Range("B2").Value = ".png"
Add-Hlink $zh "B2" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/3de2e418-5a16-46cf-9df6-9bb7cac9905f.png" ".png"

$zh.Range("D2").Value = "da2c28b8fd338c1bfb764d2362ae41a7f45681d4.png"
Add-Hlink $zh "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/HEAD/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/da2c28b8fd338c1bfb764d2362ae41a7f45681d4.png" "da2c28b8fd338c1bfb764d2362ae41a7f45681d4.png"

# Row 3: new .md dependency row
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "2016-03-14 08:42:59"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("I3").Value = "Include"

$zh.Range("A3").Value = "b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"
Add-Hlink $zh "A3" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/b1b1466b-8fdb-4993-b7a1-5da72f4af826.md" "b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"

$zh.Range("B3").Value = ".md"
Add-Hlink $zh "B3" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/b1b1466b-8fdb-4993-b7a1-5da72f4af826.md" ".md"

$zh.Range("D3").Value = "b1b1466b-8fdb-4993-b7a1-5da72f4af826.9db41ebf3b57a5038e262c17aac34d25f875dfcc.zh-cn.xlf"
Add-Hlink $zh "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/HEAD/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b1b1466b-8fdb-4993-b7a1-5da72f4af826.9db41ebf3b57a5038e262c17aac34d25f875dfcc.zh-cn.xlf" "b1b1466b-8fdb-4993-b7a1-5da72f4af826.9db41ebf3b57a5038e262c17aac34d25f875dfcc.zh-cn.xlf"

# Row 4: new .png row
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("E4").Value = "2016-03-14 08:42:59"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("I4").Value = "IsDependency"
$zh.Range("J4").Value = "e2e\b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"

$zh.Range("A4").Value = "e84bb446-0bb7-4b54-8a62-0c8b45e01771.png"
Add-Hlink $zh "A4" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/e84bb446-0bb7-4b54-8a62-0c8b45e01771.png" "e84bb446-0bb7-4b54-8a62-0c8b45e01771.png"

$zh.Range("B4").Value = ".png"
Add-Hlink $zh "B4" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/e84bb446-0bb7-4b54-8a62-0c8b45e01771.png" ".png"

$zh.Range("D4").Value = "761a28054670fd4d245a1f3e468cba17dd920655.png"
Add-Hlink $zh "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/HEAD/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/761a28054670fd4d245a1f3e468cba17dd920655.png" "761a28054670fd4d245a1f3e468cba17dd920655.png"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2 refresh
$de.Range("E2").Value = "2016-03-14 08:43:03"
$de.Range("I2").Value = "IsDependency"
$de.Range("J2").Value = "e2e\b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"

$de.Range("A2").Value = "3de2e418-5a16-46cf-9df6-9bb7cac9905f.png"
Add-Hlink $de "A2" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/3de2e418-5a16-46cf-9df6-9bb7cac9905f.png" "3de2e418-5a16-46cf-9df6-9bb7cac9905f.png"

$de.Range("B2").Value = ".png"
Add-Hlink $de "B2" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/3de2e418-5a16-46cf-9df6-9bb7cac9905f.png" ".png"

$de.Range("D2").Value = "da2c28b8fd338c1bfb764d2362ae41a7f45681d4.png"
Add-Hlink $de "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/HEAD/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/da2c28b8fd338c1bfb764d2362ae41a7f45681d4.png" "da2c28b8fd338c1bfb764d2362ae41a7f45681d4.png"

# Row 3: new .md dependency row
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "2016-03-14 08:43:03"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("I3").Value = "Include"

$de.Range("A3").Value = "b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"
Add-Hlink $de "A3" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/b1b1466b-8fdb-4993-b7a1-5da72f4af826.md" "b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"

$de.Range("B3").Value = ".md"
Add-Hlink $de "B3" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/b1b1466b-8fdb-4993-b7a1-5da72f4af826.md" ".md"

$de.Range("D3").Value = "b1b1466b-8fdb-4993-b7a1-5da72f4af826.9db41ebf3b57a5038e262c17aac34d25f875dfcc.de-de.xlf"
Add-Hlink $de "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/HEAD/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b1b1466b-8fdb-4993-b7a1-5da72f4af826.9db41ebf3b57a5038e262c17aac34d25f875dfcc.de-de.xlf" "b1b1466b-8fdb-4993-b7a1-5da72f4af826.9db41ebf3b57a5038e262c17aac34d25f875dfcc.de-de.xlf"

# Row 4: new .png row
$de.Range("C4").Value = "Ready for handoff"
$de.Range("E4").Value = "2016-03-14 08:43:03"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("I4").Value = "IsDependency"
$de.Range("J4").Value = "e2e\b1b1466b-8fdb-4993-b7a1-5da72f4af826.md"

$de.Range("A4").Value = "e84bb446-0bb7-4b54-8a62-0c8b45e01771.png"
Add-Hlink $de "A4" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/e84bb446-0bb7-4b54-8a62-0c8b45e01771.png" "e84bb446-0bb7-4b54-8a62-0c8b45e01771.png"

$de.Range("B4").Value = ".png"
Add-Hlink $de "B4" "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/e84bb446-0bb7-4b54-8a62-0c8b45e01771.png" ".png"

$de.Range("D4").Value = "761a28054670fd4d245a1f3e468cba17dd920655.png"
Add-Hlink $de "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/HEAD/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/761a28054670fd4d245a1f3e468cba17dd920655.png" "761a28054670fd4d245a1f3e468cba17dd920655.png"

"Report generated"
